$d = $word.ActiveDocument
$d.Content.Find.Execute("ceux la", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ceux-là", 2)
